# Scheduled-runner price/profit refresh for the Gungnir_Profits leve tables.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H:N) on the affected rows of each crafter sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 423.43478
$ws.Range("I19").Value = 234.5
$ws.Range("J19").Value = 490.11765
$ws.Range("K19").Value = 234.5
$ws.Range("L19").Value = 490.11765
$ws.Range("M19").Value = -59.5
$ws.Range("N19").Value = -840.11765

$ws.Range("H40").Value = 2978108
$ws.Range("I40").Value = 5683390
$ws.Range("J40").Value = 2298
$ws.Range("K40").Value = 5683390
$ws.Range("L40").Value = 2298
$ws.Range("M40").Value = -5683215
$ws.Range("N40").Value = -2648

$ws.Range("H86").Value = 15756.714
$ws.Range("I86").Value = 33667.668
$ws.Range("J86").Value = 2323.5
$ws.Range("K86").Value = 33667.668
$ws.Range("L86").Value = 2323.5
$ws.Range("M86").Value = -32544.668
$ws.Range("N86").Value = -4569.5

$ws.Range("H88").Value = 16596262
$ws.Range("I88").Value = 3000
$ws.Range("J88").Value = 18255588
$ws.Range("K88").Value = 3000
$ws.Range("L88").Value = 18255588
$ws.Range("M88").Value = -2594
$ws.Range("N88").Value = -18256400

$ws.Range("H89").Value = 15756.714
$ws.Range("I89").Value = 33667.668
$ws.Range("J89").Value = 2323.5
$ws.Range("K89").Value = 168338.34
$ws.Range("L89").Value = 11617.5
$ws.Range("M89").Value = -162722.34
$ws.Range("N89").Value = -22849.5

$ws.Range("H91").Value = 16596262
$ws.Range("I91").Value = 3000
$ws.Range("J91").Value = 18255588
$ws.Range("K91").Value = 3000
$ws.Range("L91").Value = 18255588
$ws.Range("M91").Value = -1596
$ws.Range("N91").Value = -18258396

$ws.Range("H138").Value = 1898.23
$ws.Range("I138").Value = 943.119
$ws.Range("J138").Value = 2589.862
$ws.Range("K138").Value = 2829.357
$ws.Range("L138").Value = 7769.586
$ws.Range("M138").Value = 2310.643
$ws.Range("N138").Value = -18049.586

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 50000
$ws.Range("J23").Value = 50000
$ws.Range("L23").Value = 50000
$ws.Range("N23").Value = -50518

$ws.Range("H43").Value = 7418
$ws.Range("I43").Value = 4000
$ws.Range("J43").Value = 8101.6
$ws.Range("K43").Value = 4000
$ws.Range("L43").Value = 8101.6
$ws.Range("M43").Value = -3687
$ws.Range("N43").Value = -8727.6

$ws.Range("H61").Value = 1427.0857
$ws.Range("I61").Value = 1384.7273
$ws.Range("J61").Value = 1498.7693
$ws.Range("K61").Value = 1384.7273
$ws.Range("L61").Value = 1498.7693
$ws.Range("M61").Value = -1172.7273
$ws.Range("N61").Value = -1922.7693

$ws.Range("H136").Value = 1427.0857
$ws.Range("I136").Value = 1384.7273
$ws.Range("J136").Value = 1498.7693
$ws.Range("K136").Value = 4154.1819
$ws.Range("L136").Value = 4496.3079
$ws.Range("M136").Value = -1604.1819
$ws.Range("N136").Value = -9596.3079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1164245.1
$ws.Range("I86").Value = 1235.4546
$ws.Range("J86").Value = 2585701.2
$ws.Range("K86").Value = 1235.4546
$ws.Range("L86").Value = 2585701.2
$ws.Range("M86").Value = -112.4546
$ws.Range("N86").Value = -2587947.2

$ws.Range("H89").Value = 1164245.1
$ws.Range("I89").Value = 1235.4546
$ws.Range("J89").Value = 2585701.2
$ws.Range("K89").Value = 6177.273
$ws.Range("L89").Value = 12928506
$ws.Range("M89").Value = -561.2730000000001
$ws.Range("N89").Value = -12939738

$ws.Range("H100").Value = 90000
$ws.Range("J100").Value = 90000
$ws.Range("L100").Value = 90000
$ws.Range("N100").Value = -92164

$ws.Range("H134").Value = 1567686
$ws.Range("I134").Value = 1019.881
$ws.Range("J134").Value = 3836650.8
$ws.Range("K134").Value = 3059.643
$ws.Range("L134").Value = 11509952.4
$ws.Range("M134").Value = -524.643
$ws.Range("N134").Value = -11515022.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 25220.25
$ws.Range("I22").Value = 14537.429
$ws.Range("K22").Value = 14537.429
$ws.Range("M22").Value = -14187.429

$ws.Range("H122").Value = 11364675
$ws.Range("I122").Value = 11905398
$ws.Range("J122").Value = 9500
$ws.Range("K122").Value = 35716194
$ws.Range("L122").Value = 28500
$ws.Range("M122").Value = -35713744
$ws.Range("N122").Value = -33400

$ws.Range("H132").Value = 13895096
$ws.Range("I132").Value = 1359.3636
$ws.Range("J132").Value = 25651334
$ws.Range("K132").Value = 4078.0908
$ws.Range("L132").Value = 76954002
$ws.Range("M132").Value = -1548.0908
$ws.Range("N132").Value = -76959062

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 16368.8
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 17965.334
$ws.Range("K80").Value = 6000
$ws.Range("L80").Value = 53896.00199999999
$ws.Range("M80").Value = -5064
$ws.Range("N80").Value = -55768.00199999999

$ws.Range("H83").Value = 16368.8
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 17965.334
$ws.Range("K83").Value = 18000
$ws.Range("L83").Value = 161688.006
$ws.Range("M83").Value = -13320
$ws.Range("N83").Value = -171048.006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3928.2104
$ws.Range("I126").Value = 3986.3333
$ws.Range("J126").Value = 3828.5715
$ws.Range("K126").Value = 11958.9999
$ws.Range("L126").Value = 11485.7145
$ws.Range("M126").Value = -9488.999899999999
$ws.Range("N126").Value = -16425.7145

$ws.Range("H131").Value = 49324
$ws.Range("J131").Value = 49324
$ws.Range("L131").Value = 49324
$ws.Range("N131").Value = -59404

$ws.Range("H138").Value = 57358.168
$ws.Range("J138").Value = 57358.168
$ws.Range("L138").Value = 57358.168
$ws.Range("N138").Value = -67638.168

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 576820.44
$ws.Range("I22").Value = 1151064.4
$ws.Range("J22").Value = 2576.4546
$ws.Range("K22").Value = 1151064.4
$ws.Range("L22").Value = 2576.4546
$ws.Range("M22").Value = -1150769.4
$ws.Range("N22").Value = -3166.4546

$ws.Range("H27").Value = 576820.44
$ws.Range("I27").Value = 1151064.4
$ws.Range("J27").Value = 2576.4546
$ws.Range("K27").Value = 1151064.4
$ws.Range("L27").Value = 2576.4546
$ws.Range("M27").Value = -1150957.4
$ws.Range("N27").Value = -2790.4546

$ws.Range("H46").Value = 4168140.5
$ws.Range("I46").Value = 5952946.5
$ws.Range("J46").Value = 3593.3333
$ws.Range("K46").Value = 5952946.5
$ws.Range("L46").Value = 3593.3333
$ws.Range("M46").Value = -5952758.5
$ws.Range("N46").Value = -3969.3333

$ws.Range("H123").Value = 27787.666
$ws.Range("I123").Value = 8000
$ws.Range("J123").Value = 28777.05
$ws.Range("K123").Value = 8000
$ws.Range("L123").Value = 28777.05
$ws.Range("M123").Value = -3100
$ws.Range("N123").Value = -38577.05

$ws.Range("H132").Value = 18188402
$ws.Range("I132").Value = 37038870
$ws.Range("J132").Value = 11164.75
$ws.Range("K132").Value = 111116610
$ws.Range("L132").Value = 33494.25
$ws.Range("M132").Value = -111114080
$ws.Range("N132").Value = -38554.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2648.0527
$ws.Range("I122").Value = 1782.5454
$ws.Range("J122").Value = 3838.125
$ws.Range("K122").Value = 5347.6362
$ws.Range("L122").Value = 11514.375
$ws.Range("M122").Value = -2897.6362
$ws.Range("N122").Value = -16414.375

$ws.Range("H132").Value = 22750.076
$ws.Range("I132").Value = 22995.213
$ws.Range("J132").Value = 20829.834
$ws.Range("K132").Value = 68985.639
$ws.Range("L132").Value = 62489.50199999999
$ws.Range("M132").Value = -66455.639
$ws.Range("N132").Value = -67549.502
